# Update "想去人数" (F column) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - rows 5,7,8,9,11,12,13
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 146
$wsExhibit.Range("F7").Value = 168
$wsExhibit.Range("F8").Value = 354
$wsExhibit.Range("F9").Value = 465
$wsExhibit.Range("F11").Value = 145
$wsExhibit.Range("F12").Value = 11954
$wsExhibit.Range("F13").Value = 5431

# Sheet "全部类型" (sheet4) - rows 4,7,9,10,11,13,14,15,16
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8
$wsAll.Range("F7").Value = 146
$wsAll.Range("F9").Value = 168
$wsAll.Range("F10").Value = 354
$wsAll.Range("F11").Value = 465
$wsAll.Range("F13").Value = 145
$wsAll.Range("F14").Value = 11954
$wsAll.Range("F15").Value = 8
$wsAll.Range("F16").Value = 5431
